$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defects")
$ws.Activate()

# --- Row 20: defect status changed from Open to Closed ---
$ws.Range("G2").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = "Closed"

# --- Row 21: defect status changed from Open to Closed ---
$ws.Range("G2").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("G21").Value = "Closed"

# --- New row 22: new defect entry ---
$ws.Range("H21").Copy()
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("B22").Value = "The accordion on order history tab overlaps the order details with other accordion tabs"
$ws.Range("F22").Value = "Swapnil"
$ws.Range("G22").Value = "Open"

# --- Update the view: scroll position and selection ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F26").Select() | Out-Null
